# TimeSlots.xlsx edit: "Added more users and time slots"
#
# 1) Break the external links to Schedule.xlsm (both the
#    "C:\ScheduleMeWorkspace\..." and "C:\workspace\..." sources). This drops
#    xl/externalLinks/*.xml + the <externalReferences> element.
# 2) Point the now-dangling defined names (Dates / Functions) at #REF!
# 3) Clear the empty, style-only placeholder cells in D12:J14 and replace
#    rows 12-14 with real data in columns A/B (new time slots / assignments),
#    continuing the same date style used by the rows above them.
# 4) New shared strings ("Reader1", "Acolyte EMHC MAS Reader1 Reader2") come
#    along for free as a side effect of writing those text values.
# 5) Shrink the two data-validation ranges that used to cover the now
#    filled-in rows 12-14.
# 6) Misc cosmetic bits: column A width nudge, and move the active
#    selection to B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Break external links -> removes externalLinks parts + <externalReferences>
$sources = $wb.LinkSources(1)
foreach ($src in $sources) {
    $wb.BreakLink($src, 1)
}

# --- 2) Re-point the defined names that used to resolve through those links
foreach ($n in $wb.Names) {
    $n.RefersTo = "=#REF!"
}

# --- 3) Rebuild rows 12-14 with real data instead of empty styled cells
$ws.Range("D12:J14").Clear()

$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 45300.416666666664
$ws.Range("B12").Value = "Reader1"

$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 45307.416666666664
$ws.Range("B13").Value = "Acolyte EMHC MAS Reader1 Reader2"

$ws.Range("A11").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 45276.416666666664
$ws.Range("B14").Value = "Acolyte EMHC MAS Reader1 Reader2"

# --- 5) Shrink validations that used to include rows 12-14
$ws.Range("J12:J14").Validation.Delete()
$ws.Range("E12:E14").Validation.Delete()

# --- 6) Cosmetic: widen column A slightly, move selection to B14
$ws.Columns("A:A").ColumnWidth = 25.45
[void]$ws.Range("B14").Select()
